# Update "想去人数" (column F) counts that changed between crawls.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 633
$ws1.Range("F4").Value = 663
$ws1.Range("F5").Value = 571
$ws1.Range("F7").Value = 2813
$ws1.Range("F9").Value = 7940
$ws1.Range("F12").Value = 44
$ws1.Range("F13").Value = 381

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 633
$ws4.Range("F4").Value = 663
$ws4.Range("F5").Value = 571
$ws4.Range("F9").Value = 2813
$ws4.Range("F11").Value = 7940
$ws4.Range("F14").Value = 44
$ws4.Range("F17").Value = 381
